# Remove the 4 columns (vol min, vol max, Expected part, Left/right) that sit
# between the "HU" column and the "Mandatory" column on every check-protocol
# sheet (Clinical Structures, opt structures, couch_structures). Deleting the
# entire columns C:F shifts the old "Mandatory" column (G) left into the new
# column C, exactly like selecting columns C:F in Excel and choosing Delete.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Clinical Structures", "opt structures", "couch_structures")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    [void]$ws.Activate()
    $cols = $ws.Range("C1:F1").EntireColumn
    [void]$cols.Select()
    [void]$cols.Delete()
}

# "couch_structures" (the last sheet edited) stays the active sheet/tab,
# matching the saved state in the workbook.
$ws = $wb.Worksheets.Item("couch_structures")
[void]$ws.Activate()
